# --------------------------------------------------------------------------
# Edit: rename header labels on the two existing sheets, add a new
# "PO Forecast" worksheet after "Monthly Trend" populated with forecast data.
# --------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# 1) Rename the "Requested quantity" header on the "Weekly Quantity" sheet.
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# 2) Rename the "Requested quantity" header on the "Monthly Trend" sheet.
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 3) Add a new "PO Forecast" sheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the bold/centered/bordered header style used on the other sheets.
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Match the date-style formatting used for column A on the other sheets.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A38").PasteSpecial(-4122)

# Forecast data (37 rows): ds, PO_Forecast, yhat_lower, yhat_upper
$data = New-Object 'object[,]' 37,4
$data[0,0] = 45333.99999999999
$data[0,1] = 83
$data[0,2] = -57.8503500420362
$data[0,3] = 229.8193728373158
$data[1,0] = 45340.99999999999
$data[1,1] = 87
$data[1,2] = -76.73043335459052
$data[1,3] = 238.1602943421486
$data[2,0] = 45347.99999999999
$data[2,1] = 92
$data[2,2] = -55.98969859335689
$data[2,3] = 241.477707413254
$data[3,0] = 45354.99999999999
$data[3,1] = 97
$data[3,2] = -53.03501730813685
$data[3,3] = 242.282918252319
$data[4,0] = 45361.99999999999
$data[4,1] = 102
$data[4,2] = -52.82390326439315
$data[4,3] = 255.8412007490252
$data[5,0] = 45368.99999999999
$data[5,1] = 106
$data[5,2] = -50.52851590944476
$data[5,3] = 251.9529800733644
$data[6,0] = 45375.99999999999
$data[6,1] = 111
$data[6,2] = -34.84048933277323
$data[6,3] = 271.1578146405756
$data[7,0] = 45382.99999999999
$data[7,1] = 116
$data[7,2] = -36.46348492793628
$data[7,3] = 266.5558589158263
$data[8,0] = 45389.99999999999
$data[8,1] = 121
$data[8,2] = -22.57010152578441
$data[8,3] = 278.308087956552
$data[9,0] = 45396.99999999999
$data[9,1] = 126
$data[9,2] = -17.87219633684695
$data[9,3] = 276.3483917915082
$data[10,0] = 45403.99999999999
$data[10,1] = 130
$data[10,2] = -15.79441137727736
$data[10,3] = 290.0162643395922
$data[11,0] = 45410.99999999999
$data[11,1] = 135
$data[11,2] = -14.47703430235295
$data[11,3] = 279.2973192723124
$data[12,0] = 45424.99999999999
$data[12,1] = 145
$data[12,2] = -3.802590525145463
$data[12,3] = 301.0364788457757
$data[13,0] = 45431.99999999999
$data[13,1] = 149
$data[13,2] = 2.929783554903956
$data[13,3] = 308.0541494137693
$data[14,0] = 45438.99999999999
$data[14,1] = 154
$data[14,2] = -3.90117821742051
$data[14,3] = 296.5432706996621
$data[15,0] = 45445.99999999999
$data[15,1] = 159
$data[15,2] = 12.05465942428991
$data[15,3] = 314.7479398482934
$data[16,0] = 45459.99999999999
$data[16,1] = 168
$data[16,2] = 15.6068752121735
$data[16,3] = 311.9072864254423
$data[17,0] = 45466.99999999999
$data[17,1] = 173
$data[17,2] = 15.81722404546971
$data[17,3] = 334.9684607765215
$data[18,0] = 45473.99999999999
$data[18,1] = 178
$data[18,2] = 23.94067736529591
$data[18,3] = 329.5781816968621
$data[19,0] = 45487.99999999999
$data[19,1] = 187
$data[19,2] = 35.74306072965542
$data[19,3] = 336.819558042649
$data[20,0] = 45529.99999999999
$data[20,1] = 216
$data[20,2] = 55.30936791308882
$data[20,3] = 372.2464968449867
$data[21,0] = 45543.99999999999
$data[21,1] = 225
$data[21,2] = 66.42102622274307
$data[21,3] = 373.7035878391399
$data[22,0] = 45550.99999999999
$data[22,1] = 230
$data[22,2] = 71.48441564006464
$data[22,3] = 368.8644674392694
$data[23,0] = 45557.99999999999
$data[23,1] = 235
$data[23,2] = 82.39041033170271
$data[23,3] = 384.5193883100656
$data[24,0] = 45564.99999999999
$data[24,1] = 240
$data[24,2] = 89.4315049632621
$data[24,3] = 382.0744780932127
$data[25,0] = 45592.99999999999
$data[25,1] = 259
$data[25,2] = 114.9587201826617
$data[25,3] = 412.0842798306267
$data[26,0] = 45599.99999999999
$data[26,1] = 263
$data[26,2] = 110.1306978658112
$data[26,3] = 419.0063207556536
$data[27,0] = 45613.99999999999
$data[27,1] = 273
$data[27,2] = 117.5199754934683
$data[27,3] = 435.2741188017335
$data[28,0] = 45620.99999999999
$data[28,1] = 278
$data[28,2] = 128.9268553349866
$data[28,3] = 436.5780523963947
$data[29,0] = 45627.99999999999
$data[29,1] = 282
$data[29,2] = 142.6771976846175
$data[29,3] = 447.0092805206103
$data[30,0] = 45634.99999999999
$data[30,1] = 287
$data[30,2] = 138.9093070874919
$data[30,3] = 433.3624191419881
$data[31,0] = 45641.99999999999
$data[31,1] = 292
$data[31,2] = 139.4397881255803
$data[31,3] = 441.9839127992196
$data[32,0] = 45648.99999999999
$data[32,1] = 297
$data[32,2] = 144.603066493347
$data[32,3] = 451.5030164769115
$data[33,0] = 45655.99999999999
$data[33,1] = 301
$data[33,2] = 149.9593847438544
$data[33,3] = 445.4380889145255
$data[34,0] = 45662.99999999999
$data[34,1] = 306
$data[34,2] = 165.293754904533
$data[34,3] = 461.802037832822
$data[35,0] = 45669.99999999999
$data[35,1] = 311
$data[35,2] = 157.1265760264202
$data[35,3] = 461.5951076398223
$data[36,0] = 45676.99999999999
$data[36,1] = 316
$data[36,2] = 167.2081416509683
$data[36,3] = 475.4618270792108

$wsForecast.Range("A2:D38").Value2 = $data

Write-Output "PO Forecast sheet created with $($wsForecast.UsedRange.Rows.Count) rows"
